$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.889.78"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.877.92"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.14"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4599"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3878"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07876"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9846"
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.75"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.873.10"
$ws.Range("E12").Value = "  -2.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.992"
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.661"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06949"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.20"
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009974"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.95"
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.903.52"
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.252"
$ws.Range("E22").Value = "  -1.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.95"
$ws.Range("E23").Value = "  -0.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.089"
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.85"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.31"
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.995"
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.927"
$ws.Range("E28").Value = "  -1.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.35"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09346"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9040"
$ws.Range("E31").Value = "  -2.74%  "
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.319"
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.263"
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.187"
$ws.Range("E35").Value = "  +2.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05767"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02072"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.001"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.675"
$ws.Range("E39").Value = "  -1.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5654"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1766"
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.652"
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.258"
$ws.Range("E43").Value = "  +1.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.93"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5350"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07033"
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "113.18"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.508"
$ws.Range("E49").Value = "  +1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.75"
$ws.Range("E51").Value = "  -0.29%  "
